$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "reapers"

$headers = @("author", "creator", "producer", "subject", "title")
$values = @(
    "A. Murugan; S.Anu H. Nair; K. P. Sanal Kumar",
    "Arbortext Advanced Print Publisher 9.1.440/W Unicode",
    "Acrobat Distiller 9.0.0 (Windows); modified using iText® 5.3.5 ©2000-2012 1T3XT BVBA (AGPL-version)",
    "J Med Syst, doi:10.1007/s10916-019-1400-8",
    "Detection of Skin Cancer Using SVM, Random Forest and kNN Classifiers"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $headers[$i]
    $ws.Cells.Item(2, $col).Value = $values[$i]
}

# Build the header format (bold, boxed border, centered/top-aligned) on a
# scratch cell so the whole combination becomes a single new style record,
# then copy+paste-special (formats only) onto the real header row so the
# target cells pick up that one combined style instead of accumulating a
# fresh style per individual property assignment.
$scratch = $ws.Range("Z100")
$scratch.Value = "tmp"
$scratch.Font.Bold = $true
$scratch.Borders.LineStyle = 1
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4160

$scratch.Copy()
$ws.Range("A1:E1").PasteSpecial(-4122)
$scratch.Clear()

$excel.CutCopyMode = $false
